# Edit script for "Algorithm - Medium - 2 Key Keyboards.docx"
#
# Implements the diff:
#   - mini_step(B) = mini_step(A) + D   ->   mini_step(B) = min(mini_step(A) + D, B)
#   - adds a new explanatory sentence + line break after that formula
#   - merges "(A)" + " <" into "(A) <" around "Then as A < B -> mini_step(A) <"
#   - merges "(A) + D" + " is much better than B" into one run
#   - appends a new "min(mini_step(A) + D, B) = mini_step(A) + D" line
#
# All inserts are done with Find (locate only) + Collapse + InsertBefore/
# InsertAfter so that existing <w:proofErr> spellcheck markers around the
# "mini_step" occurrences are never split apart (which the naive
# Find.Execute(..., Replace:=wdReplaceOne) text-replace would do for a
# match that starts/ends inside a proofErr-wrapped run).

$d = $word.ActiveDocument
$wdCollapseStart = 1
$wdCollapseEnd = 0

# 1) "(B) = " -> "(B) = min("
$r1 = $d.Content
[void]$r1.Find.Execute("(B) = ")
$r1.Collapse($wdCollapseEnd)
$r1.InsertAfter("min(")

# 2) first "(A) + D" found after that point -> "(A) + D, B)"
$r2 = $d.Range($r1.End, $d.Content.End)
[void]$r2.Find.Execute("(A) + D")
$r2.Collapse($wdCollapseEnd)
$r2.InsertAfter(", B)")

# 3) Insert the new parenthetical sentence (plus its trailing line break)
#    right before "Then as A < B -> "
$r3 = $d.Content
[void]$r3.Find.Execute("Then as A " + [char]60 + " B -" + [char]62)
$r3.Collapse($wdCollapseStart)
$newSentence = "(As we can attain B, we can either do it B times or go to its factor first then make the copy)" + [char]11
$r3.InsertBefore($newSentence)

# 4) Merge "(A)" + " <" (around "mini_step(A) <= B/2") into a single run "(A) <"
[void]$d.Content.Find.Execute("(A) " + [char]60, $true, $false, $false, $false, $false, $true, 1, $false, "(A) " + [char]60, 2)

# 5) Merge "(A) + D" + " is much better than B" into a single run
[void]$d.Content.Find.Execute("(A) + D is much better than B", $true, $false, $false, $false, $false, $true, 1, $false, "(A) + D is much better than B", 2)

# 6) Append the new closing line: a line break followed by
#    "min(mini_step(A) + D, B) = mini_step(A) + D"
$r6 = $d.Content
[void]$r6.Find.Execute("is much better than B")
$r6.Collapse($wdCollapseEnd)
$r6.InsertAfter([char]11 + "min(mini_step(A) + D, B) = mini_step(A) + D")
